# Refresh the "cryptos" price/volume table (Price = column D, Volume(1h) = column E)
# with the latest scraped values from the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Both columns hold plain text (prices like "28.114.24" / "1.005" are not meant to be
# numbers, and the percentage strings keep their padding spaces), so force each write
# to stay text - a leading apostrophe stops Excel's automatic number detection - and
# then reset the cell style so it doesn't keep the "stored as text" quote-prefix style.
function Set-TextValue {
    param($ws, $cellRef, $text)
    $r = $ws.Range($cellRef)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "28.114.24"
Set-TextValue $ws "E2" "  +2.91%  "
Set-TextValue $ws "D3" "1.776.13"
Set-TextValue $ws "E3" "  -0.48%  "
Set-TextValue $ws "D4" "1.005"
Set-TextValue $ws "E4" "  +0.17%  "
Set-TextValue $ws "D5" "339.03"
Set-TextValue $ws "E5" "  -0.32%  "
Set-TextValue $ws "E6" "  +0.39%  "
Set-TextValue $ws "D7" "0.3818"
Set-TextValue $ws "E7" "  -3.54%  "
Set-TextValue $ws "D8" "0.3416"
Set-TextValue $ws "E8" "  -1.49%  "
Set-TextValue $ws "D9" "46.96"
Set-TextValue $ws "E9" "  -2.25%  "
Set-TextValue $ws "D10" "1.144"
Set-TextValue $ws "E10" "  -4.31%  "
Set-TextValue $ws "D11" "0.07378"
Set-TextValue $ws "E11" "  -1.30%  "
Set-TextValue $ws "D12" "23.31"
Set-TextValue $ws "E12" "  +7.19%  "
Set-TextValue $ws "D13" "1.005"
Set-TextValue $ws "E13" "  +0.41%  "
Set-TextValue $ws "D14" "6.387"
Set-TextValue $ws "E14" "  -1.32%  "
Set-TextValue $ws "D15" "7.298"
Set-TextValue $ws "E15" "  +2.66%  "
Set-TextValue $ws "D16" "1.777.44"
Set-TextValue $ws "E16" "  -0.34%  "
Set-TextValue $ws "D17" "0.00001077"
Set-TextValue $ws "E17" "  -1.66%  "
Set-TextValue $ws "D18" "0.06660"
Set-TextValue $ws "E18" "  -0.47%  "
Set-TextValue $ws "D19" "82.31"
Set-TextValue $ws "E19" "  -2.78%  "
Set-TextValue $ws "D20" "1.003"
Set-TextValue $ws "D21" "17.35"
Set-TextValue $ws "E21" "  -2.55%  "
Set-TextValue $ws "D22" "6.391"
Set-TextValue $ws "E22" "  -1.74%  "
Set-TextValue $ws "D23" "28.124.22"
Set-TextValue $ws "E23" "  +2.89%  "
Set-TextValue $ws "D24" "12.09"
Set-TextValue $ws "E24" "  -2.40%  "
Set-TextValue $ws "D25" "2.392"
Set-TextValue $ws "E25" "  -0.03%  "
Set-TextValue $ws "D26" "20.66"
Set-TextValue $ws "E26" "  -2.73%  "
Set-TextValue $ws "D27" "1.424"
Set-TextValue $ws "E27" "  -2.30%  "
Set-TextValue $ws "D28" "2.396"
Set-TextValue $ws "E28" "  -3.69%  "
Set-TextValue $ws "D29" "153.86"
Set-TextValue $ws "E29" "  -2.51%  "
Set-TextValue $ws "D30" "1.980.30"
Set-TextValue $ws "E30" "  -0.23%  "
Set-TextValue $ws "D31" "134.85"
Set-TextValue $ws "E31" "  -0.90%  "
Set-TextValue $ws "D32" "4.025"
Set-TextValue $ws "E32" "  +0.04%  "
Set-TextValue $ws "D33" "6.059"
Set-TextValue $ws "E33" "  +1.55%  "
Set-TextValue $ws "D34" "0.08861"
Set-TextValue $ws "E34" "  +0.40%  "
Set-TextValue $ws "D35" "12.70"
Set-TextValue $ws "E35" "  -2.08%  "
Set-TextValue $ws "D36" "0.02405"
Set-TextValue $ws "E36" "  -0.59%  "
Set-TextValue $ws "D37" "0.6824"
Set-TextValue $ws "E37" "  -0.01%  "
Set-TextValue $ws "D38" "5.317"
Set-TextValue $ws "E38" "  -1.58%  "
Set-TextValue $ws "D39" "0.06348"
Set-TextValue $ws "E39" "  -1.97%  "
Set-TextValue $ws "D40" "0.2153"
Set-TextValue $ws "E40" "  -2.57%  "
Set-TextValue $ws "D41" "1.240"
Set-TextValue $ws "E41" "  -1.06%  "
Set-TextValue $ws "D42" "1.499"
Set-TextValue $ws "E42" "  -7.31%  "
Set-TextValue $ws "D43" "8.218"
Set-TextValue $ws "E43" "  -1.33%  "
Set-TextValue $ws "D44" "1.003"
Set-TextValue $ws "E44" "  +0.39%  "
Set-TextValue $ws "D45" "14.06"
Set-TextValue $ws "E45" "  -2.94%  "
Set-TextValue $ws "D46" "0.6243"
Set-TextValue $ws "E46" "  -2.10%  "
Set-TextValue $ws "D47" "3.866"
Set-TextValue $ws "E47" "  -0.24%  "
Set-TextValue $ws "D48" "132.67"
Set-TextValue $ws "E48" "  +0.60%  "
Set-TextValue $ws "D49" "2.061"
Set-TextValue $ws "E49" "  -3.41%  "
Set-TextValue $ws "D50" "0.07498"
Set-TextValue $ws "E50" "  +4.73%  "
Set-TextValue $ws "D51" "1.201"
Set-TextValue $ws "E51" "  +4.26%  "
